$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing shared-string values to their upper-case forms
$ws.Range("A2").Value = "FOO"
$ws.Range("A3").Value = "BAR"

# Add new row with the concatenated value
$ws.Range("A4").Value = "FOOBAR"

# Move selection to the newly added cell, matching the authored selection
$ws.Range("A4").Select()
